$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 830.61536
$ws.Range("I33").Value = 862.3333
$ws.Range("K33").Value = 862.3333
$ws.Range("M33").Value = -633.3333
$ws.Range("H41").Value = 719.6
$ws.Range("I41").Value = 431.14285
$ws.Range("K41").Value = 431.14285
$ws.Range("M41").Value = 8.85714999999999
$ws.Range("H53").Value = 363.78946
$ws.Range("I53").Value = 308.41666
$ws.Range("K53").Value = 308.41666
$ws.Range("M53").Value = 328.58334
$ws.Range("H95").Value = 23097.2
$ws.Range("J95").Value = 23097.2
$ws.Range("L95").Value = 23097.2
$ws.Range("N95").Value = -28589.2
$ws.Range("H112").Value = 1122.7826
$ws.Range("J112").Value = 1108.2858
$ws.Range("L112").Value = 3324.8574
$ws.Range("N112").Value = -5540.857400000001
$ws.Range("H131").Value = 8822.477000000001
$ws.Range("I131").Value = 2052.7
$ws.Range("J131").Value = 14976.818
$ws.Range("K131").Value = 6158.099999999999
$ws.Range("L131").Value = 44930.454
$ws.Range("M131").Value = -1118.099999999999
$ws.Range("N131").Value = -55010.454
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2781.8572
$ws.Range("I2").Value = 1549.6666
$ws.Range("J2").Value = 4999.8
$ws.Range("K2").Value = 1549.6666
$ws.Range("L2").Value = 4999.8
$ws.Range("M2").Value = -1436.6666
$ws.Range("N2").Value = -5225.8
$ws.Range("H32").Value = 3402.976
$ws.Range("I32").Value = 3168.9023
$ws.Range("J32").Value = 13000
$ws.Range("K32").Value = 3168.9023
$ws.Range("L32").Value = 13000
$ws.Range("M32").Value = -2881.9023
$ws.Range("N32").Value = -13574
$ws.Range("H61").Value = 6336.875
$ws.Range("J61").Value = 7999
$ws.Range("L61").Value = 7999
$ws.Range("N61").Value = -8423
$ws.Range("H74").Value = 2750.7646
$ws.Range("I74").Value = 2622.6875
$ws.Range("K74").Value = 2622.6875
$ws.Range("M74").Value = -1748.6875
$ws.Range("H77").Value = 2750.7646
$ws.Range("I77").Value = 2622.6875
$ws.Range("K77").Value = 13113.4375
$ws.Range("M77").Value = -8745.4375
$ws.Range("H113").Value = 139964
$ws.Range("J113").Value = 139964
$ws.Range("L113").Value = 139964
$ws.Range("N113").Value = -148642
$ws.Range("H116").Value = 2781.8572
$ws.Range("I116").Value = 1549.6666
$ws.Range("J116").Value = 4999.8
$ws.Range("K116").Value = 1549.6666
$ws.Range("L116").Value = 4999.8
$ws.Range("M116").Value = 744.3334
$ws.Range("N116").Value = -9587.799999999999
$ws.Range("H122").Value = 2699.75
$ws.Range("I122").Value = 1899.5
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 5698.5
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -3248.5
$ws.Range("N122").Value = -15400
$ws.Range("H133").Value = 500261
$ws.Range("J133").Value = 500261
$ws.Range("L133").Value = 500261
$ws.Range("N133").Value = -505321
$ws.Range("H136").Value = 6336.875
$ws.Range("J136").Value = 7999
$ws.Range("L136").Value = 23997
$ws.Range("N136").Value = -29097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2781.8572
$ws.Range("I3").Value = 1549.6666
$ws.Range("J3").Value = 4999.8
$ws.Range("K3").Value = 1549.6666
$ws.Range("L3").Value = 4999.8
$ws.Range("M3").Value = -1435.6666
$ws.Range("N3").Value = -5227.8
$ws.Range("H96").Value = 12969
$ws.Range("I96").Value = 12969
$ws.Range("K96").Value = 12969
$ws.Range("M96").Value = -10223
$ws.Range("H105").Value = 3080.875
$ws.Range("I105").Value = 3110.75
$ws.Range("J105").Value = 2991.25
$ws.Range("K105").Value = 3110.75
$ws.Range("L105").Value = 2991.25
$ws.Range("M105").Value = -1363.75
$ws.Range("N105").Value = -6485.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 163795.28
$ws.Range("I94").Value = 370319
$ws.Range("K94").Value = 370319
$ws.Range("M94").Value = -369868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 3387.625
$ws.Range("I38").Value = 2849.8333
$ws.Range("K38").Value = 8549.499899999999
$ws.Range("M38").Value = -8202.499899999999
$ws.Range("H69").Value = 2663.8333
$ws.Range("I69").Value = 1991.5
$ws.Range("K69").Value = 5974.5
$ws.Range("M69").Value = -5163.5
$ws.Range("H70").Value = 11547.5
$ws.Range("I70").Value = 2460
$ws.Range("J70").Value = 17000
$ws.Range("K70").Value = 7380
$ws.Range("L70").Value = 51000
$ws.Range("M70").Value = -7065
$ws.Range("N70").Value = -51630
$ws.Range("H72").Value = 2663.8333
$ws.Range("I72").Value = 1991.5
$ws.Range("K72").Value = 17923.5
$ws.Range("M72").Value = -13867.5
$ws.Range("H73").Value = 11547.5
$ws.Range("I73").Value = 2460
$ws.Range("J73").Value = 17000
$ws.Range("K73").Value = 7380
$ws.Range("L73").Value = 51000
$ws.Range("M73").Value = -6288
$ws.Range("N73").Value = -53184
$ws.Range("H129").Value = 2533.2307
$ws.Range("I129").Value = 708
$ws.Range("K129").Value = 2124
$ws.Range("M129").Value = 2876
$ws.Range("H134").Value = 16169.091
$ws.Range("I134").Value = 6620
$ws.Range("J134").Value = 19750
$ws.Range("K134").Value = 19860
$ws.Range("L134").Value = 59250
$ws.Range("M134").Value = -14790
$ws.Range("N134").Value = -69390

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5731.75
$ws.Range("I70").Value = 5488
$ws.Range("J70").Value = 5766.5713
$ws.Range("K70").Value = 5488
$ws.Range("L70").Value = 5766.5713
$ws.Range("M70").Value = -5218
$ws.Range("N70").Value = -6306.5713
$ws.Range("H73").Value = 5731.75
$ws.Range("I73").Value = 5488
$ws.Range("J73").Value = 5766.5713
$ws.Range("K73").Value = 5488
$ws.Range("L73").Value = 5766.5713
$ws.Range("M73").Value = -4552
$ws.Range("N73").Value = -7638.5713
$ws.Range("H102").Value = 2354.2222
$ws.Range("I102").Value = 2091.7334
$ws.Range("K102").Value = 2091.7334
$ws.Range("M102").Value = -469.7334000000001
$ws.Range("H107").Value = 434.4
$ws.Range("I107").Value = 390.125
$ws.Range("J107").Value = 611.5
$ws.Range("K107").Value = 390.125
$ws.Range("L107").Value = 611.5
$ws.Range("M107").Value = 1529.875
$ws.Range("N107").Value = -4451.5
$ws.Range("I113").Value = 3298.5
$ws.Range("J113").Value = 3999.5
$ws.Range("K113").Value = 3298.5
$ws.Range("L113").Value = 3999.5
$ws.Range("M113").Value = -1128.5
$ws.Range("N113").Value = -8339.5
$ws.Range("H122").Value = 1000
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1138.8
$ws.Range("I16").Value = 995.8
$ws.Range("K16").Value = 995.8
$ws.Range("M16").Value = -825.8
$ws.Range("H82").Value = 1628
$ws.Range("I82").Value = 1599.4
$ws.Range("K82").Value = 1599.4
$ws.Range("M82").Value = -1238.4
$ws.Range("H85").Value = 1628
$ws.Range("I85").Value = 1599.4
$ws.Range("K85").Value = 1599.4
$ws.Range("M85").Value = -351.4000000000001
$ws.Range("H128").Value = 59999
$ws.Range("J128").Value = 59999
$ws.Range("L128").Value = 59999
$ws.Range("N128").Value = -69959
$ws.Range("H136").Value = 7897.6665
$ws.Range("I136").Value = 8384.875
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 25154.625
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -22604.625
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1000000.8
$ws.Range("J5").Value = 1000000.8
$ws.Range("L5").Value = 1000000.8
$ws.Range("N5").Value = -1000224.8
$ws.Range("H62").Value = 7762
$ws.Range("I62").Value = 6900.727
$ws.Range("J62").Value = 12499
$ws.Range("K62").Value = 6900.727
$ws.Range("L62").Value = 12499
$ws.Range("M62").Value = -6276.727
$ws.Range("N62").Value = -13747
$ws.Range("H65").Value = 7762
$ws.Range("I65").Value = 6900.727
$ws.Range("J65").Value = 12499
$ws.Range("K65").Value = 34503.635
$ws.Range("L65").Value = 62495
$ws.Range("M65").Value = -31383.635
$ws.Range("N65").Value = -68735
$ws.Range("H82").Value = 38825.25
$ws.Range("I82").Value = 30000
$ws.Range("J82").Value = 41767
$ws.Range("K82").Value = 30000
$ws.Range("L82").Value = 41767
$ws.Range("M82").Value = -29617
$ws.Range("N82").Value = -42533
$ws.Range("H85").Value = 38825.25
$ws.Range("I85").Value = 30000
$ws.Range("J85").Value = 41767
$ws.Range("K85").Value = 30000
$ws.Range("L85").Value = 41767
$ws.Range("M85").Value = -28674
$ws.Range("N85").Value = -44419
$ws.Range("H132").Value = 10295.4
$ws.Range("I132").Value = 9119.5
$ws.Range("K132").Value = 27358.5
$ws.Range("M132").Value = -24828.5
$ws.Range("H136").Value = 2629.2856
$ws.Range("I136").Value = 2693.077
$ws.Range("K136").Value = 8079.231000000001
$ws.Range("M136").Value = -5529.231000000001
